$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-10-14 Monday" "2024-10-15 Tuesday"

Replace-Text "277×7=" "485×5="
Replace-Text "855×8=" "315×4="
Replace-Text "812×4=" "383×3="
Replace-Text "347×7=" "255×2="
Replace-Text "964×5=" "296×5="

Replace-Text "932×8=" "578×9="
Replace-Text "144×8=" "981×8="
Replace-Text "716×5=" "999×4="
Replace-Text "235×8=" "628×7="
Replace-Text "564×7=" "928×3="

Replace-Text "355×8=" "972×2="
Replace-Text "449×9=" "657×6="
Replace-Text "462×9=" "267×5="
Replace-Text "457×3=" "864×8="
Replace-Text "256×9=" "444×3="

Replace-Text "342×8=" "209×4="
Replace-Text "963×5=" "883×7="
Replace-Text "389×9=" "694×6="
Replace-Text "691×3=" "655×2="
Replace-Text "609×8=" "695×5="

Replace-Text "334×8=" "259×7="
Replace-Text "702×3=" "535×2="
Replace-Text "139×6=" "712×4="
Replace-Text "445×7=" "381×6="
Replace-Text "958×3=" "381×6="
